$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "76.519.45"
$ws.Range("E2").Value = "  +0.77%  "
$ws.Range("D3").Value = "3.042.87"
$ws.Range("E3").Value = "  +4.29%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'201.14"
$ws.Range("E5").Value = "  +0.99%  "
$ws.Range("D6").Value = "'629.83"
$ws.Range("E6").Value = "  +5.11%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "'0.552"
$ws.Range("E8").Value = "  +0.37%  "
$ws.Range("D9").Value = "'0.203"
$ws.Range("E9").Value = "  +2.37%  "
$ws.Range("D10").Value = "3.047.60"
$ws.Range("E10").Value = "  +4.57%  "
$ws.Range("D11").Value = "'0.435"
$ws.Range("E11").Value = "  -1.12%  "
$ws.Range("E12").Value = "  -0.25%  "
$ws.Range("D13").Value = "'5.08"
$ws.Range("E13").Value = "  +4.05%  "
$ws.Range("D14").Value = "3.607.95"
$ws.Range("E14").Value = "  +4.45%  "
$ws.Range("D15").Value = "'29.37"
$ws.Range("E15").Value = "  +6.73%  "
$ws.Range("D16").Value = "76.513.22"
$ws.Range("E16").Value = "  +0.90%  "
$ws.Range("D17").Value = "'0.0000189"
$ws.Range("E17").Value = "  -0.62%  "
$ws.Range("D18").Value = "3.041.91"
$ws.Range("E18").Value = "  +4.40%  "
$ws.Range("D19").Value = "'13.57"
$ws.Range("E19").Value = "  +6.12%  "
$ws.Range("D20").Value = "'9.04"
$ws.Range("E20").Value = "  +1.02%  "
$ws.Range("D21").Value = "'376.75"
$ws.Range("E21").Value = "  -0.50%  "
$ws.Range("D22").Value = "'4.35"
$ws.Range("E22").Value = "  +2.56%  "
$ws.Range("D23").Value = "'2.28"
$ws.Range("E23").Value = "  -1.30%  "
$ws.Range("D24").Value = "'73.28"
$ws.Range("E24").Value = "  +2.72%  "
$ws.Range("D25").Value = "3.206.45"
$ws.Range("E25").Value = "  +4.55%  "
$ws.Range("E26").Value = "  +5.07%  "
$ws.Range("E27").Value = "  -0.13%  "
$ws.Range("D28").Value = "'9.88"
$ws.Range("E28").Value = "  +1.94%  "
$ws.Range("E29").Value = "  +0.63%  "
$ws.Range("D30").Value = "'0.998"
$ws.Range("E30").Value = "  -0.21%  "
$ws.Range("D31").Value = "'8.30"
$ws.Range("E31").Value = "  +7.51%  "
$ws.Range("D32").Value = "'1.40"
$ws.Range("E32").Value = "  -1.52%  "
$ws.Range("D33").Value = "'512.02"
$ws.Range("E33").Value = "  +1.05%  "
$ws.Range("D34").Value = "'1.94"
$ws.Range("E34").Value = "  +6.75%  "
$ws.Range("E35").Value = "  +0.25%  "
$ws.Range("D36").Value = "'20.88"
$ws.Range("E36").Value = "  +3.45%  "
$ws.Range("D37").Value = "'164.27"
$ws.Range("E37").Value = "  -0.26%  "
$ws.Range("B38").Value = "WhiteBITCoin"
$ws.Range("C38").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D38").Value = "'20.01"
$ws.Range("E38").Value = "  +1.72%  "
$ws.Range("B39").Value = "PolygonEcosystemToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D39").Value = "'0.383"
$ws.Range("E39").Value = "  +11.50%  "
$ws.Range("D40").Value = "'191.38"
$ws.Range("E40").Value = "  +5.53%  "
$ws.Range("E41").Value = "  -3.61%  "
$ws.Range("D42").Value = "'0.112"
$ws.Range("E42").Value = "  -1.21%  "
$ws.Range("E43").Value = "  +0.15%  "
$ws.Range("E44").Value = "  +1.35%  "
$ws.Range("D45").Value = "'43.15"
$ws.Range("E45").Value = "  +7.27%  "
$ws.Range("D46").Value = "'1.25"
$ws.Range("E46").Value = "  +4.49%  "
$ws.Range("D47").Value = "'1.66"
$ws.Range("E47").Value = "  -0.14%  "
$ws.Range("D48").Value = "'0.611"
$ws.Range("E48").Value = "  +6.73%  "
$ws.Range("D49").Value = "'0.707"
$ws.Range("E49").Value = "  +6.61%  "
$ws.Range("E50").Value = "  +1.01%  "
$ws.Range("E51").Value = "  +4.49%  "
